$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs",  "Lrfn3", "Lrfn3", "ECs",  2, 0.6666666666666666, 0.589336,           1.768008,           0.1262869893132262, 0.1262869893132262, 2, 0.6666666666666666, 0.589336,           1.768008,           0.1262869893132262, 0.1262869893132262, 0.347316920896,     3.125852288064,      0.0159484036697989,  0.0159484036697989),
    @("ECs",  "Lrfn3", "Lrfn3", "FAPs", 2, 0.6666666666666666, 0.589336,           1.768008,           0.1262869893132262, 0.1262869893132262, 3, 1,                   3.172374,           9.517122000000001, 0.6797982160186321, 0.6797982160186321, 1.869594203664,     16.826347832976,     0.0858496700414952,  0.0858496700414952),
    @("ECs",  "Lrfn3", "Lrfn3", "sCs",  2, 0.6666666666666666, 0.589336,           1.768008,           0.1262869893132262, 0.1262869893132262, 3, 1,                   0.9049306666666667, 2.714792,           0.1939147946681417, 0.1939147946681417, 0.5333082193706666, 4.799773974336,      0.02448891560193206, 0.02448891560193206),
    @("FAPs", "Lrfn3", "Lrfn3", "ECs",  3, 1,                   3.172374,           9.517122000000001, 0.6797982160186321, 0.6797982160186321, 2, 0.6666666666666666, 0.589336,           1.768008,           0.1262869893132262, 0.1262869893132262, 1.869594203664,     16.826347832976,     0.0858496700414952,  0.0858496700414952),
    @("FAPs", "Lrfn3", "Lrfn3", "FAPs", 3, 1,                   3.172374,           9.517122000000001, 0.6797982160186321, 0.6797982160186321, 3, 1,                   3.172374,           9.517122000000001, 0.6797982160186321, 0.6797982160186321, 10.063956795876,    90.57561116288402,   0.4621256145021148,  0.4621256145021148),
    @("FAPs", "Lrfn3", "Lrfn3", "sCs",  3, 1,                   3.172374,           9.517122000000001, 0.6797982160186321, 0.6797982160186321, 3, 1,                   0.9049306666666667, 2.714792,           0.1939147946681417, 0.1939147946681417, 2.870778518736,     25.837006668624,     0.1318229314750221,  0.1318229314750221),
    @("sCs",  "Lrfn3", "Lrfn3", "ECs",  3, 1,                   0.9049306666666667, 2.714792,           0.1939147946681417, 0.1939147946681417, 2, 0.6666666666666666, 0.589336,           1.768008,           0.1262869893132262, 0.1262869893132262, 0.5333082193706666, 4.799773974336,      0.02448891560193206, 0.02448891560193206),
    @("sCs",  "Lrfn3", "Lrfn3", "FAPs", 3, 1,                   0.9049306666666667, 2.714792,           0.1939147946681417, 0.1939147946681417, 3, 1,                   3.172374,           9.517122000000001, 0.6797982160186321, 0.6797982160186321, 2.870778518736,     25.837006668624,     0.1318229314750221,  0.1318229314750221),
    @("sCs",  "Lrfn3", "Lrfn3", "sCs",  3, 1,                   0.9049306666666667, 2.714792,           0.1939147946681417, 0.1939147946681417, 3, 1,                   0.9049306666666667, 2.714792,           0.1939147946681417, 0.1939147946681417, 0.8188995114737778, 7.370095603264001,   0.03760294759118756, 0.03760294759118756)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $rowData = $rows[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowData[$j]
    }
}
